$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# Original text (subtitle placeholder on slide 1):
#   P1: "Tom Dietterich, Oregon State"
#   P2: "Chad Scherrer, Galois"
#   P3: "Roddy Collins, " + "Kitware"
#
# Target text (author order fixed):
#   P1: "Tom Dietterich, Oregon State"              (unchanged)
#   P2: "Roddy " + "Collins, " + "Kitware"          (was P3, now 2nd)
#   P3: "Chad Scherrer" + ", " + "Galois"           (was P2, now 3rd/last)
#
# Swap the paragraph contents in place (keeps paragraph count/order of marks,
# only the text itself moves), then re-split each paragraph into its runs.

# P2 currently holds "Chad Scherrer, Galois" (21 chars, starts at char 30)
$p2 = $tr.Characters(30, 21)
$p2.Text = "Roddy Collins, Kitware"

# P3 currently holds "Roddy Collins, Kitware" (22 chars, now starts at char 53)
$p3 = $tr.Characters(53, 22)
$p3.Text = "Chad Scherrer, Galois"

# Split P2 "Roddy Collins, Kitware" into "Roddy " / "Collins, " / "Kitware"
$roddyWord = $tr.Characters(30, 6)
$roddyWord.Text = "Roddy "
$collinsWord = $tr.Characters(36, 9)
$collinsWord.Text = "Collins, "
$kitwareWord = $tr.Characters(45, 7)
$kitwareWord.Text = "Kitware"

# Split P3 "Chad Scherrer, Galois" into "Chad Scherrer" / ", " / "Galois"
$chadName = $tr.Characters(53, 13)
$chadName.Text = "Chad Scherrer"
$comma = $tr.Characters(66, 2)
$comma.Text = ", "
$galois = $tr.Characters(68, 6)
$galois.Text = "Galois"
